# visits and cashflow generation
# Remove the "Vizsla" breed row from the dog sheet (row 142) and shift
# everything below it up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")

# Find the row containing "Vizsla" in column A and delete the entire row,
# shifting the rows below it upward.
$target = $ws.Range("A1:A200").Find("Vizsla")
if ($target -ne $null) {
    $target.EntireRow.Delete()
}

# Update the active selection to match the post-edit workbook state.
$ws.Range("J10").Select()

$wb.Save()
